# Daily attendance processing - 2025-12-29 11:31:55
# Normalizes the ordering of names/emails in the "Recorded By" column (G).
# Applies the following exact text substitutions wherever found:
#   "System, system, backup@backdoor.com" -> "System, backup@backdoor.com, system"
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System"         -> "System, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    switch ($val) {
        "System, system, backup@backdoor.com" { $cell.Value2 = "System, backup@backdoor.com, system" }
        "System, dnasr281@gmail.com"          { $cell.Value2 = "dnasr281@gmail.com, System" }
        "backup@backdoor.com, System"         { $cell.Value2 = "System, backup@backdoor.com" }
    }
}

$wb.Save()
